$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the stray ">:(" note that lived outside the table (column E).
$ws.Range("E29").ClearContents()

# Make room for a new "user_id" field in the "posts" table (new row 14),
# pushing everything below it down by one row.
$ws.Rows.Item(14).Insert()

# --- posts table: add user_id right after id, keep the rest as-is ---
$ws.Range("C14").Value = "user_id"
$ws.Range("C15").Value = "title"
$ws.Range("C16").Value = "description"
$ws.Range("C17").Value = "shortdesc"
$ws.Range("C18").Value = "elkeszites"
$ws.Range("C19").Value = "adag"
$ws.Range("C20").Value = "ingredients"
$ws.Range("C21").Value = "datum"
$ws.Range("C22").Value = "points"

# --- comments table (renamed from "kommentek" / "tartalom" -> "comment"), now placed where favorites/recept_kepek used to start ---
$ws.Range("B23").Value = "comments"
$ws.Range("C23").Value = "id"
$ws.Range("C24").Value = "post_id"
$ws.Range("C25").Value = "user_id"
$ws.Range("C26").Value = "comment"
$ws.Range("C27").Value = "datum"
$ws.Range("C28").Value = "points"

# --- follows table (renamed from "kovetesek" / "user" / "kovetett_user") ---
$ws.Range("B29").Value = "follows"
$ws.Range("C29").Value = "user_id"
$ws.Range("C30").Value = "kovetett_user_id"

# --- favorites table (unchanged content, shifted into new position) ---
$ws.Range("B31").Value = "favorites"
$ws.Range("C31").Value = "user_id"
$ws.Range("C32").Value = "post_id"

# --- recept_kepek table (unchanged content, shifted into new position) ---
$ws.Range("B33").Value = "recept_kepek"
$ws.Range("C33").Value = "post_id"
$ws.Range("C34").Value = "filename"

# Re-shape the table-name merges in column B to match the new row groups.
$ws.Range("B23:B24").UnMerge()
$ws.Range("B25:B26").UnMerge()
$ws.Range("B27:B28").UnMerge()
$ws.Range("B23:B28").Merge()

$ws.Range("B29:B34").UnMerge()
$ws.Range("B29:B30").Merge()
$ws.Range("B31:B32").Merge()
$ws.Range("B33:B34").Merge()

# Column C got a bit wider to fit "kovetett_user_id".
$ws.Columns.Item(3).ColumnWidth = 16.28515625

# View settings: scroll down a little and move the selection.
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("D24").Select()
